$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2911.7246
$ws.Range("I76").Value = 2732.6086
$ws.Range("J76").Value = 3269.9565
$ws.Range("K76").Value = 2732.6086
$ws.Range("L76").Value = 3269.9565
$ws.Range("M76").Value = -2417.6086
$ws.Range("N76").Value = -3899.9565
$ws.Range("H79").Value = 2911.7246
$ws.Range("I79").Value = 2732.6086
$ws.Range("J79").Value = 3269.9565
$ws.Range("K79").Value = 2732.6086
$ws.Range("L79").Value = 3269.9565
$ws.Range("M79").Value = -1640.6086
$ws.Range("N79").Value = -5453.9565
$ws.Range("H100").Value = 2180.3572
$ws.Range("I100").Value = 1501.4286
$ws.Range("J100").Value = 2859.2856
$ws.Range("K100").Value = 1501.4286
$ws.Range("L100").Value = 2859.2856
$ws.Range("M100").Value = -960.4286
$ws.Range("N100").Value = -3941.2856
$ws.Range("H113").Value = 2405.6667
$ws.Range("I113").Value = 2050.8333
$ws.Range("J113").Value = 2642.2222
$ws.Range("K113").Value = 2050.8333
$ws.Range("L113").Value = 2642.2222
$ws.Range("M113").Value = 1203.1667
$ws.Range("N113").Value = -9150.2222
$ws.Range("H116").Value = 2405.3125
$ws.Range("I116").Value = 2098.8462
$ws.Range("J116").Value = 3733.3333
$ws.Range("K116").Value = 2098.8462
$ws.Range("L116").Value = 3733.3333
$ws.Range("M116").Value = 1343.1538
$ws.Range("N116").Value = -10617.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8398.741
$ws.Range("I32").Value = 8381.598
$ws.Range("K32").Value = 8381.598
$ws.Range("M32").Value = -8094.598
$ws.Range("H110").Value = 1463.375
$ws.Range("I110").Value = 1377.25
$ws.Range("J110").Value = 1549.5
$ws.Range("K110").Value = 1377.25
$ws.Range("L110").Value = 1549.5
$ws.Range("M110").Value = 667.75
$ws.Range("N110").Value = -5639.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3377.439
$ws.Range("I105").Value = 1657.4117
$ws.Range("J105").Value = 4595.7915
$ws.Range("K105").Value = 1657.4117
$ws.Range("L105").Value = 4595.7915
$ws.Range("M105").Value = 89.58829999999989
$ws.Range("N105").Value = -8089.7915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7096236
$ws.Range("I31").Value = 4768.9653
$ws.Range("J31").Value = 18521376
$ws.Range("K31").Value = 4768.9653
$ws.Range("L31").Value = 18521376
$ws.Range("M31").Value = -4473.9653
$ws.Range("N31").Value = -18521966
$ws.Range("H34").Value = 7096236
$ws.Range("I34").Value = 4768.9653
$ws.Range("J34").Value = 18521376
$ws.Range("K34").Value = 4768.9653
$ws.Range("L34").Value = 18521376
$ws.Range("M34").Value = -4566.9653
$ws.Range("N34").Value = -18521780
$ws.Range("H105").Value = 1598.75
$ws.Range("I105").Value = 1765
$ws.Range("J105").Value = 1100
$ws.Range("K105").Value = 1765
$ws.Range("L105").Value = 1100
$ws.Range("M105").Value = -18
$ws.Range("N105").Value = -4594

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4833.3335
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 5600
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 16800
$ws.Range("M76").Value = -2617
$ws.Range("N76").Value = -17566
$ws.Range("H79").Value = 4833.3335
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 5600
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 16800
$ws.Range("M79").Value = -1674
$ws.Range("N79").Value = -19452
$ws.Range("H88").Value = 2823.3333
$ws.Range("J88").Value = 2823.3333
$ws.Range("L88").Value = 8469.999899999999
$ws.Range("N88").Value = -9325.999899999999
$ws.Range("H91").Value = 2823.3333
$ws.Range("J91").Value = 2823.3333
$ws.Range("L91").Value = 8469.999899999999
$ws.Range("N91").Value = -11433.9999
$ws.Range("H94").Value = 7484.8
$ws.Range("I94").Value = 1024
$ws.Range("J94").Value = 9100
$ws.Range("K94").Value = 3072
$ws.Range("L94").Value = 27300
$ws.Range("M94").Value = -2396
$ws.Range("N94").Value = -28652
$ws.Range("H100").Value = 9200
$ws.Range("J100").Value = 9200
$ws.Range("L100").Value = 27600
$ws.Range("N100").Value = -29222
$ws.Range("H103").Value = 4098.423
$ws.Range("I103").Value = 693.3333
$ws.Range("J103").Value = 5901.1177
$ws.Range("K103").Value = 2079.9999
$ws.Range("L103").Value = 17703.3531
$ws.Range("M103").Value = -1200.9999
$ws.Range("N103").Value = -19461.3531
$ws.Range("H109").Value = 3882.5454
$ws.Range("I109").Value = 1063.2
$ws.Range("J109").Value = 4711.7646
$ws.Range("K109").Value = 3189.6
$ws.Range("L109").Value = 14135.2938
$ws.Range("M109").Value = -2149.6
$ws.Range("N109").Value = -16215.2938
$ws.Range("H112").Value = 4312.9487
$ws.Range("I112").Value = 3335
$ws.Range("J112").Value = 4456.7646
$ws.Range("K112").Value = 10005
$ws.Range("L112").Value = 13370.2938
$ws.Range("M112").Value = -8897
$ws.Range("N112").Value = -15586.2938
$ws.Range("H113").Value = 785.4286
$ws.Range("I113").Value = 512.93335
$ws.Range("J113").Value = 1215.6842
$ws.Range("K113").Value = 1538.80005
$ws.Range("L113").Value = 3647.0526
$ws.Range("M113").Value = 631.1999499999999
$ws.Range("N113").Value = -7987.0526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13780
$ws.Range("I70").Value = 200000
$ws.Range("J70").Value = 3978.9473
$ws.Range("K70").Value = 200000
$ws.Range("L70").Value = 3978.9473
$ws.Range("M70").Value = -199730
$ws.Range("N70").Value = -4518.9473
$ws.Range("H73").Value = 13780
$ws.Range("I73").Value = 200000
$ws.Range("J73").Value = 3978.9473
$ws.Range("K73").Value = 200000
$ws.Range("L73").Value = 3978.9473
$ws.Range("M73").Value = -199064
$ws.Range("N73").Value = -5850.9473
$ws.Range("H113").Value = 334303.66
$ws.Range("I113").Value = 500605.5
$ws.Range("K113").Value = 500605.5
$ws.Range("M113").Value = -498435.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8001.5557
$ws.Range("I122").Value = 9152.799999999999
$ws.Range("J122").Value = 6562.5
$ws.Range("K122").Value = 27458.4
$ws.Range("L122").Value = 19687.5
$ws.Range("M122").Value = -25008.4
$ws.Range("N122").Value = -24587.5
